# RPA datasets push 2024-07-26
# Insert a new record at the top of the data table (row 2), pushing all
# existing rows down by one. The new row describes the SK증권제13호스팩
# offering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 2 (everything currently at row 2+ shifts
# down to row 3+). EntireRow.Insert() mimics Excel's native "Insert Row"
# and inherits the header row's formatting, so explicitly clear it back to
# the plain/default look used by the rest of the data rows.
$ws.Range("A2").EntireRow.Insert()
$ws.Range("A2:T2").ClearFormats()

# Columns A, D and E hold date-shaped text (e.g. "2024-07-15") that must
# stay literal text, matching every other row in the sheet (no numeric
# date serials). Writing them as a formula that evaluates to the literal
# string, then collapsing the formula to its cached value via a
# copy/paste-values round trip, avoids Excel's automatic "this looks like
# a date" coercion that a plain .Value assignment would trigger.
$ws.Range("A2").Formula = '="2024-07-15"'
$ws.Range("D2").Formula = '="2024-07-18"'
$ws.Range("E2").Formula = '="2024-07-25"'
$ws.Range("A2:E2").Copy()
$ws.Range("A2:E2").PasteSpecial(-4163)

$ws.Range("B2").Value = "SK증권제13호스팩"
$ws.Range("C2").Value = "SK"
$ws.Range("F2").Value = 8000000
$ws.Range("G2").Value = 4000000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "1654.28 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"
